$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 417, pushing existing rows 417-454 down to 418-455
$ws.Rows.Item(417).Insert()

# Populate the newly inserted row 417 with the new data record
$ws.Range("A417").Value = 4
$ws.Range("B417").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C417").Value = "Los Lagos"
$ws.Range("D417").Value = 45166
$ws.Range("E417").Value = 10
$ws.Range("F417").Value = "Fruta"
$ws.Range("G417").Value = 100108
$ws.Range("H417").Value = "Tropicales y subtropicales"
$ws.Range("I417").Value = 100108005
$ws.Range("J417").Value = "Piña"
$ws.Range("K417").Value = "Caramelo"
$ws.Range("L417").Value = "Segunda"
$ws.Range("M417").Value = 50
$ws.Range("N417").Value = 25000
$ws.Range("O417").Value = 25000
$ws.Range("P417").Value = 25000
$ws.Range("Q417").Value = "$/caja 14 unidades"
$ws.Range("R417").Value = "Ecuador"
$ws.Range("S417").Value = 1786
$ws.Range("T417").Value = 14
